$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential disclaimer date in the shared text cell (A16)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03071948499600446
$ws.Range("E2").Value = -0.002044293015332266

$ws.Range("D3").Value = 0.0234910595790896
$ws.Range("E3").Value = -0.0007058823529412228

$ws.Range("D4").Value = 0.05038082565047302
$ws.Range("E4").Value = -0.0002380952380951484

$ws.Range("D5").Value = 0.1381175517317821
$ws.Range("E5").Value = -0.003379465722561981

$ws.Range("D6").Value = 0.02876505235948436
$ws.Range("E6").Value = 0.03516819571865448

$ws.Range("D7").Value = 0.1212911321397092
$ws.Range("E7").Value = -0.001354512991011081

$ws.Range("D8").Value = 0.100970042577913
$ws.Range("E8").Value = -0.002576370997423694

$ws.Range("D9").Value = 0.0279433650839945
$ws.Range("E9").Value = 0.002146383344065139

$ws.Range("D10").Value = 0.1212009311750175
$ws.Range("E10").Value = 0.003104987385988744

$ws.Range("D11").Value = 0.2551734722138104
$ws.Range("E11").Value = -0.007267950963222414

$ws.Range("D12").Value = 0.1019470824927219
$ws.Range("E12").Value = 0.00670930232558109

$ws.Range("E13").Value = -0.0007052421813581455

# Restore worksheet protection (sheet was protected before this edit; the
# original password hash cannot be reconstructed, so re-protect without one)
$ws.Protect()
